$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2160193333333333
$ws.Range("H2").Value = 0.648058
$ws.Range("I2").Value = 0.02486881244588016
$ws.Range("J2").Value = 0.02486881244588016
$ws.Range("M2").Value = 128.4548946666667
$ws.Range("N2").Value = 385.364684
$ws.Range("O2").Value = 0.2815548034715028
$ws.Range("P2").Value = 0.2815548034715028
$ws.Range("Q2").Value = 27.74874070929689
$ws.Range("R2").Value = 249.738666383672
$ws.Range("S2").Value = 0.00700193360076945
$ws.Range("T2").Value = 0.007001933600769453
$ws.Range("G3").Value = 0.2160193333333333
$ws.Range("H3").Value = 0.648058
$ws.Range("I3").Value = 0.02486881244588016
$ws.Range("J3").Value = 0.02486881244588016
$ws.Range("N3").Value = 457.183265
$ws.Range("O3").Value = 0.3340268313936494
$ws.Range("P3").Value = 0.3340268313936494
$ws.Range("Q3").Value = 32.92014137215222
$ws.Range("R3").Value = 296.28127234937
$ws.Range("S3").Value = 0.008306850621820301
$ws.Range("T3").Value = 0.008306850621820303
$ws.Range("G4").Value = 0.2160193333333333
$ws.Range("H4").Value = 0.648058
$ws.Range("I4").Value = 0.02486881244588016
$ws.Range("J4").Value = 0.02486881244588016
$ws.Range("M4").Value = 70.798157
$ws.Range("N4").Value = 212.394471
$ws.Range("O4").Value = 0.1551794599342134
$ws.Range("P4").Value = 0.1551794599342134
$ws.Range("Q4").Value = 15.29377067636867
$ws.Range("R4").Value = 137.643936087318
$ws.Range("S4").Value = 0.003859128884556927
$ws.Range("T4").Value = 0.003859128884556928
$ws.Range("G5").Value = 0.2160193333333333
$ws.Range("H5").Value = 0.648058
$ws.Range("I5").Value = 0.02486881244588016
$ws.Range("J5").Value = 0.02486881244588016
$ws.Range("M5").Value = 20.703408
$ws.Range("N5").Value = 62.110224
$ws.Range("O5").Value = 0.04537891674549766
$ws.Range("P5").Value = 0.04537891674549767
$ws.Range("Q5").Value = 4.472336393888
$ws.Range("R5").Value = 40.25102754499201
$ws.Range("S5").Value = 0.001128519769540992
$ws.Range("T5").Value = 0.001128519769540992
$ws.Range("G6").Value = 0.2160193333333333
$ws.Range("H6").Value = 0.648058
$ws.Range("I6").Value = 0.02486881244588016
$ws.Range("J6").Value = 0.02486881244588016
$ws.Range("M6").Value = 83.88319133333333
$ws.Range("N6").Value = 251.649574
$ws.Range("O6").Value = 0.1838599884551367
$ws.Range("P6").Value = 0.1838599884551367
$ws.Range("Q6").Value = 18.12039106969911
$ws.Range("R6").Value = 163.083519627292
$ws.Range("S6").Value = 0.004572379569192484
$ws.Range("T6").Value = 0.004572379569192486
$ws.Range("I7").Value = 0.9551554900377276
$ws.Range("J7").Value = 0.9551554900377278
$ws.Range("M7").Value = 128.4548946666667
$ws.Range("N7").Value = 385.364684
$ws.Range("O7").Value = 0.2815548034715028
$ws.Range("P7").Value = 0.2815548034715028
$ws.Range("Q7").Value = 1065.767096349995
$ws.Range("R7").Value = 9591.903867149955
$ws.Range("S7").Value = 0.2689286162822994
$ws.Range("T7").Value = 0.2689286162822994
$ws.Range("I8").Value = 0.9551554900377276
$ws.Range("J8").Value = 0.9551554900377278
$ws.Range("N8").Value = 457.183265
$ws.Range("O8").Value = 0.3340268313936494
$ws.Range("P8").Value = 0.3340268313936494
$ws.Range("S8").Value = 0.3190475618255506
$ws.Range("T8").Value = 0.3190475618255507
$ws.Range("I9").Value = 0.9551554900377276
$ws.Range("J9").Value = 0.9551554900377278
$ws.Range("M9").Value = 70.798157
$ws.Range("N9").Value = 212.394471
$ws.Range("O9").Value = 0.1551794599342134
$ws.Range("P9").Value = 0.1551794599342134
$ws.Range("Q9").Value = 587.3995413613543
$ws.Range("R9").Value = 5286.595872252188
$ws.Range("S9").Value = 0.1482205130972535
$ws.Range("T9").Value = 0.1482205130972535
$ws.Range("I10").Value = 0.9551554900377276
$ws.Range("J10").Value = 0.9551554900377278
$ws.Range("M10").Value = 20.703408
$ws.Range("N10").Value = 62.110224
$ws.Range("O10").Value = 0.04537891674549766
$ws.Range("P10").Value = 0.04537891674549767
$ws.Range("Q10").Value = 171.772442661424
$ws.Range("R10").Value = 1545.951983952816
$ws.Range("S10").Value = 0.04334392146142706
$ws.Range("T10").Value = 0.04334392146142708
$ws.Range("I11").Value = 0.9551554900377276
$ws.Range("J11").Value = 0.9551554900377278
$ws.Range("M11").Value = 83.88319133333333
$ws.Range("N11").Value = 251.649574
$ws.Range("O11").Value = 0.1838599884551367
$ws.Range("P11").Value = 0.1838599884551367
$ws.Range("Q11").Value = 695.9637115571628
$ws.Range("R11").Value = 6263.673404014465
$ws.Range("S11").Value = 0.175614877371197
$ws.Range("T11").Value = 0.175614877371197
$ws.Range("G12").Value = 0.173491
$ws.Range("H12").Value = 0.520473
$ws.Range("I12").Value = 0.01997281943922393
$ws.Range("J12").Value = 0.01997281943922393
$ws.Range("M12").Value = 128.4548946666667
$ws.Range("N12").Value = 385.364684
$ws.Range("O12").Value = 0.2815548034715028
$ws.Range("P12").Value = 0.2815548034715028
$ws.Range("Q12").Value = 22.28576813061466
$ws.Range("R12").Value = 200.571913175532
$ws.Range("S12").Value = 0.005623443251982504
$ws.Range("T12").Value = 0.005623443251982506
$ws.Range("G13").Value = 0.173491
$ws.Range("H13").Value = 0.520473
$ws.Range("I13").Value = 0.01997281943922393
$ws.Range("J13").Value = 0.01997281943922393
$ws.Range("N13").Value = 457.183265
$ws.Range("O13").Value = 0.3340268313936494
$ws.Range("P13").Value = 0.3340268313936494
$ws.Range("Q13").Value = 26.43906060937166
$ws.Range("R13").Value = 237.951545484345
$ws.Range("S13").Value = 0.006671457591281454
$ws.Range("T13").Value = 0.006671457591281457
$ws.Range("G14").Value = 0.173491
$ws.Range("H14").Value = 0.520473
$ws.Range("I14").Value = 0.01997281943922393
$ws.Range("J14").Value = 0.01997281943922393
$ws.Range("M14").Value = 70.798157
$ws.Range("N14").Value = 212.394471
$ws.Range("O14").Value = 0.1551794599342134
$ws.Range("P14").Value = 0.1551794599342134
$ws.Range("Q14").Value = 12.282843056087
$ws.Range("R14").Value = 110.545587504783
$ws.Range("S14").Value = 0.003099371333942328
$ws.Range("T14").Value = 0.003099371333942329
$ws.Range("G15").Value = 0.173491
$ws.Range("H15").Value = 0.520473
$ws.Range("I15").Value = 0.01997281943922393
$ws.Range("J15").Value = 0.01997281943922393
$ws.Range("M15").Value = 20.703408
$ws.Range("N15").Value = 62.110224
$ws.Range("O15").Value = 0.04537891674549766
$ws.Range("P15").Value = 0.04537891674549767
$ws.Range("Q15").Value = 3.591854957328
$ws.Range("R15").Value = 32.326694615952
$ws.Range("S15").Value = 0.0009063449105053999
$ws.Range("T15").Value = 0.0009063449105054004
$ws.Range("G16").Value = 0.173491
$ws.Range("H16").Value = 0.520473
$ws.Range("I16").Value = 0.01997281943922393
$ws.Range("J16").Value = 0.01997281943922393
$ws.Range("M16").Value = 83.88319133333333
$ws.Range("N16").Value = 251.649574
$ws.Range("O16").Value = 0.1838599884551367
$ws.Range("P16").Value = 0.1838599884551367
$ws.Range("Q16").Value = 14.55297874761133
$ws.Range("R16").Value = 130.976808728502
$ws.Range("S16").Value = 0.003672202351512241
$ws.Range("T16").Value = 0.003672202351512242
$ws.Range("G17").Value = 0.000025
$ws.Range("H17").Value = 0.00007499999999999999
$ws.Range("I17").Value = 0.000002878077168156263
$ws.Range("J17").Value = 0.000002878077168156264
$ws.Range("M17").Value = 128.4548946666667
$ws.Range("N17").Value = 385.364684
$ws.Range("O17").Value = 0.2815548034715028
$ws.Range("P17").Value = 0.2815548034715028
$ws.Range("Q17").Value = 0.003211372366666666
$ws.Range("R17").Value = 0.0289023513
$ws.Range("S17").Value = 0.000000810336451456056
$ws.Range("T17").Value = 0.0000008103364514560562
$ws.Range("G18").Value = 0.000025
$ws.Range("H18").Value = 0.00007499999999999999
$ws.Range("I18").Value = 0.000002878077168156263
$ws.Range("J18").Value = 0.000002878077168156264
$ws.Range("N18").Value = 457.183265
$ws.Range("O18").Value = 0.3340268313936494
$ws.Range("P18").Value = 0.3340268313936494
$ws.Range("Q18").Value = 0.003809860541666666
$ws.Range("R18").Value = 0.034288744875
$ws.Range("S18").Value = 0.000000961354996985644
$ws.Range("T18").Value = 0.0000009613549969856442
$ws.Range("G19").Value = 0.000025
$ws.Range("H19").Value = 0.00007499999999999999
$ws.Range("I19").Value = 0.000002878077168156263
$ws.Range("J19").Value = 0.000002878077168156264
$ws.Range("M19").Value = 70.798157
$ws.Range("N19").Value = 212.394471
$ws.Range("O19").Value = 0.1551794599342134
$ws.Range("P19").Value = 0.1551794599342134
$ws.Range("Q19").Value = 0.001769953925
$ws.Range("R19").Value = 0.015929585325
$ws.Range("S19").Value = 0.0000004466184606034791
$ws.Range("T19").Value = 0.0000004466184606034793
$ws.Range("G20").Value = 0.000025
$ws.Range("H20").Value = 0.00007499999999999999
$ws.Range("I20").Value = 0.000002878077168156263
$ws.Range("J20").Value = 0.000002878077168156264
$ws.Range("M20").Value = 20.703408
$ws.Range("N20").Value = 62.110224
$ws.Range("O20").Value = 0.04537891674549766
$ws.Range("P20").Value = 0.04537891674549767
$ws.Range("Q20").Value = 0.0005175851999999999
$ws.Range("R20").Value = 0.0046582668
$ws.Range("S20").Value = 0.0000001306040242008807
$ws.Range("T20").Value = 0.0000001306040242008808
$ws.Range("G21").Value = 0.000025
$ws.Range("H21").Value = 0.00007499999999999999
$ws.Range("I21").Value = 0.000002878077168156263
$ws.Range("J21").Value = 0.000002878077168156264
$ws.Range("M21").Value = 83.88319133333333
$ws.Range("N21").Value = 251.649574
$ws.Range("O21").Value = 0.1838599884551367
$ws.Range("P21").Value = 0.1838599884551367
$ws.Range("Q21").Value = 0.002097079783333333
$ws.Range("R21").Value = 0.01887371805
$ws.Range("S21").Value = 0.0000005291632349102032
$ws.Range("T21").Value = 0.0000005291632349102032
